$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.466.57'
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.019.88'
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.67'
$ws.Range("E5").Value = '  +1.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.62'
$ws.Range("E6").Value = '  +3.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.020.00'
$ws.Range("E8").Value = '  +0.78%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.520'
$ws.Range("E9").Value = '  -1.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.42'
$ws.Range("E10").Value = '  +11.27%  '
$ws.Range("E11").Value = '  +2.22%  '
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000233'
$ws.Range("E13").Value = '  +2.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.60'
$ws.Range("E14").Value = '  +0.51%  '
$ws.Range("E15").Value = '  +2.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.519.45'
$ws.Range("E16").Value = '  +0.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.03'
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.335.05'
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.020.81'
$ws.Range("E19").Value = '  +0.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '449.26'
$ws.Range("E20").Value = '  -1.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.25'
$ws.Range("E21").Value = '  +2.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.692'
$ws.Range("E22").Value = '  +0.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.44'
$ws.Range("E23").Value = '  +0.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.41'
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("E25").Value = '  +2.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.83'
$ws.Range("E26").Value = '  +11.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.04'
$ws.Range("E27").Value = '  -0.96%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("E29").Value = '  +2.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.19'
$ws.Range("E31").Value = '  +3.89%  '
$ws.Range("E32").Value = '  +2.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.59'
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("E34").Value = '  +1.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0853'
$ws.Range("E35").Value = '  +6.45%  '
$ws.Range("E36").Value = '  +0.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.86'
$ws.Range("E37").Value = '  +2.22%  '
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.04'
$ws.Range("E38").Value = '  +6.01%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.08'
$ws.Range("E39").Value = '  -1.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.15'
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("E41").Value = '  -1.56%  '
$ws.Range("E42").Value = '  +0.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.34'
$ws.Range("E43").Value = '  +11.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.285'
$ws.Range("E44").Value = '  +7.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '393.56'
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("E46").Value = '  -0.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.746.23'
$ws.Range("E47").Value = '  +0.58%  '
$ws.Range("E48").Value = '  +4.00%  '
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.108'
$ws.Range("E51").Value = '  -0.85%  '
